$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 385-457: shift data down by 2 rows (columns D,H,I,J,K,L,M,N,P,Q only);
# rows 385-386 receive brand-new entries for a new reporting date.
$updates = @(
    @{ Row=385; D=45005; H="Sin especificar"; I="Primera"; J=130; K=4000; L=5000; M=4500; N="`$/caja 70 unidades"; P=64; Q=70 },
    @{ Row=386; D=45005; H="Sin especificar"; I="Segunda"; J=150; K=3000; L=4000; M=3500; N="`$/caja 100 unidades"; P=35; Q=100 },
    @{ Row=387; D=44529; H="Sin especificar"; I="Primera"; J=120; K=4000; L=4500; M=4250; N="`$/caja 70 unidades"; P=61; Q=70 },
    @{ Row=388; D=44529; H="Sin especificar"; I="Segunda"; J=120; K=3500; L=4000; M=3750; N="`$/caja 100 unidades"; P=38; Q=100 },
    @{ Row=389; D=44411; H="Alaska"; I="Segunda"; J=160; K=7000; L=8000; M=7500; N="`$/caja 100 unidades"; P=75; Q=100 },
    @{ Row=390; D=44603; H="Sin especificar"; I="Primera"; J=130; K=10000; L=11000; M=10500; N="`$/caja 70 unidades"; P=150; Q=70 },
    @{ Row=391; D=44603; H="Sin especificar"; I="Segunda"; J=160; K=8000; L=9000; M=8500; N="`$/caja 100 unidades"; P=85; Q=100 },
    @{ Row=392; D=44650; H="Sin especificar"; I="Primera"; J=120; K=12000; L=13000; M=12500; N="`$/caja 70 unidades"; P=179; Q=70 },
    @{ Row=393; D=44650; H="Sin especificar"; I="Segunda"; J=140; K=9000; L=10000; M=9500; N="`$/caja 100 unidades"; P=95; Q=100 },
    @{ Row=394; D=44169; H="Sin especificar"; I="Primera"; J=160; K=3500; L=4000; M=3750; N="`$/caja 70 unidades"; P=54; Q=70 },
    @{ Row=395; D=44169; H="Sin especificar"; I="Segunda"; J=160; K=3000; L=3500; M=3250; N="`$/caja 100 unidades"; P=32; Q=100 },
    @{ Row=396; D=44165; H="Sin especificar"; I="Primera"; J=160; K=3000; L=3500; M=3250; N="`$/caja 70 unidades"; P=46; Q=70 },
    @{ Row=397; D=44165; H="Sin especificar"; I="Segunda"; J=200; K=2500; L=3000; M=2750; N="`$/caja 100 unidades"; P=28; Q=100 },
    @{ Row=398; D=44701; H="Sin especificar"; I="Primera"; J=160; K=13000; L=14000; M=13500; N="`$/caja 70 unidades"; P=193; Q=70 },
    @{ Row=399; D=44701; H="Sin especificar"; I="Segunda"; J=160; K=10000; L=11000; M=10500; N="`$/caja 100 unidades"; P=105; Q=100 },
    @{ Row=400; D=44550; H="Sin especificar"; I="Primera"; J=120; K=4000; L=4500; M=4250; N="`$/caja 70 unidades"; P=61; Q=70 },
    @{ Row=401; D=44550; H="Sin especificar"; I="Segunda"; J=160; K=3000; L=3500; M=3250; N="`$/caja 100 unidades"; P=32; Q=100 },
    @{ Row=402; D=44797; H="Sin especificar"; I="Segunda"; J=150; K=13000; L=14000; M=13500; N="`$/caja 100 unidades"; P=135; Q=100 },
    @{ Row=403; D=44204; H="Sin especificar"; I="Primera"; J=120; K=6000; L=7000; M=6500; N="`$/caja 70 unidades"; P=93; Q=70 },
    @{ Row=404; D=44204; H="Sin especificar"; I="Segunda"; J=120; K=5000; L=6000; M=5500; N="`$/caja 100 unidades"; P=55; Q=100 },
    @{ Row=405; D=44572; H="Sin especificar"; I="Primera"; J=120; K=4000; L=4500; M=4250; N="`$/caja 70 unidades"; P=61; Q=70 },
    @{ Row=406; D=44572; H="Sin especificar"; I="Segunda"; J=120; K=3000; L=3500; M=3250; N="`$/caja 100 unidades"; P=32; Q=100 },
    @{ Row=407; D=44760; H="Sin especificar"; I="Primera"; J=150; K=14000; L=15000; M=14500; N="`$/caja 70 unidades"; P=207; Q=70 },
    @{ Row=408; D=44760; H="Sin especificar"; I="Segunda"; J=170; K=11000; L=12000; M=11500; N="`$/caja 100 unidades"; P=115; Q=100 },
    @{ Row=409; D=44291; H="Sin especificar"; I="Primera"; J=120; K=11000; L=12000; M=11500; N="`$/caja 70 unidades"; P=164; Q=70 },
    @{ Row=410; D=44291; H="Sin especificar"; I="Segunda"; J=160; K=9000; L=10000; M=9500; N="`$/caja 100 unidades"; P=95; Q=100 },
    @{ Row=411; D=44876; H="Sin especificar"; I="Primera"; J=150; K=15000; L=16000; M=15500; N="`$/caja 70 unidades"; P=221; Q=70 },
    @{ Row=412; D=44876; H="Sin especificar"; I="Segunda"; J=160; K=11000; L=12000; M=11500; N="`$/caja 100 unidades"; P=115; Q=100 },
    @{ Row=413; D=44386; H="Sin especificar"; I="Primera"; J=130; K=13000; L=14000; M=13500; N="`$/caja 70 unidades"; P=193; Q=70 },
    @{ Row=414; D=44386; H="Sin especificar"; I="Segunda"; J=160; K=10000; L=11000; M=10500; N="`$/caja 100 unidades"; P=105; Q=100 },
    @{ Row=415; D=44428; H="Sin especificar"; I="Primera"; J=120; K=11000; L=12000; M=11500; N="`$/caja 70 unidades"; P=164; Q=70 },
    @{ Row=416; D=44428; H="Sin especificar"; I="Segunda"; J=120; K=8000; L=9000; M=8500; N="`$/caja 100 unidades"; P=85; Q=100 },
    @{ Row=417; D=44253; H="Sin especificar"; I="Primera"; J=120; K=7000; L=8000; M=7500; N="`$/caja 70 unidades"; P=107; Q=70 },
    @{ Row=418; D=44253; H="Sin especificar"; I="Segunda"; J=120; K=6000; L=7000; M=6500; N="`$/caja 100 unidades"; P=65; Q=100 },
    @{ Row=419; D=44998; H="Sin especificar"; I="Primera"; J=160; K=6000; L=7000; M=6500; N="`$/caja 70 unidades"; P=93; Q=70 },
    @{ Row=420; D=44998; H="Sin especificar"; I="Segunda"; J=170; K=4000; L=5000; M=4500; N="`$/caja 100 unidades"; P=45; Q=100 },
    @{ Row=421; D=44740; H="Sin especificar"; I="Primera"; J=130; K=14000; L=15000; M=14500; N="`$/caja 70 unidades"; P=207; Q=70 },
    @{ Row=422; D=44740; H="Sin especificar"; I="Segunda"; J=150; K=12000; L=13000; M=12500; N="`$/caja 100 unidades"; P=125; Q=100 },
    @{ Row=423; D=44176; H="Sin especificar"; I="Primera"; J=120; K=6000; L=6500; M=6250; N="`$/caja 70 unidades"; P=89; Q=70 },
    @{ Row=424; D=44176; H="Sin especificar"; I="Segunda"; J=160; K=5000; L=5500; M=5250; N="`$/caja 100 unidades"; P=52; Q=100 },
    @{ Row=425; D=44795; H="Sin especificar"; I="Primera"; J=150; K=18000; L=19000; M=18500; N="`$/caja 70 unidades"; P=264; Q=70 },
    @{ Row=426; D=44795; H="Sin especificar"; I="Segunda"; J=200; K=14000; L=15000; M=14500; N="`$/caja 100 unidades"; P=145; Q=100 },
    @{ Row=427; D=44651; H="Sin especificar"; I="Primera"; J=120; K=11000; L=12000; M=11500; N="`$/caja 70 unidades"; P=164; Q=70 },
    @{ Row=428; D=44651; H="Sin especificar"; I="Segunda"; J=130; K=9000; L=10000; M=9500; N="`$/caja 100 unidades"; P=95; Q=100 },
    @{ Row=429; D=44610; H="Sin especificar"; I="Primera"; J=120; K=8000; L=9000; M=8500; N="`$/caja 70 unidades"; P=121; Q=70 },
    @{ Row=430; D=44610; H="Sin especificar"; I="Segunda"; J=160; K=6000; L=7000; M=6500; N="`$/caja 100 unidades"; P=65; Q=100 },
    @{ Row=431; D=44939; H="Sin especificar"; I="Primera"; J=150; K=10000; L=11000; M=10500; N="`$/caja 70 unidades"; P=150; Q=70 },
    @{ Row=432; D=44939; H="Sin especificar"; I="Segunda"; J=160; K=7000; L=8000; M=7500; N="`$/caja 100 unidades"; P=75; Q=100 },
    @{ Row=433; D=44568; H="Sin especificar"; I="Primera"; J=160; K=4000; L=4500; M=4250; N="`$/caja 70 unidades"; P=61; Q=70 },
    @{ Row=434; D=44568; H="Sin especificar"; I="Segunda"; J=160; K=3000; L=3500; M=3250; N="`$/caja 100 unidades"; P=32; Q=100 },
    @{ Row=435; D=44473; H="Sin especificar"; I="Primera"; J=130; K=14000; L=15000; M=14500; N="`$/caja 70 unidades"; P=207; Q=70 },
    @{ Row=436; D=44473; H="Sin especificar"; I="Segunda"; J=120; K=11000; L=12000; M=11500; N="`$/caja 100 unidades"; P=115; Q=100 },
    @{ Row=437; D=44620; H="Sin especificar"; I="Primera"; J=130; K=11000; L=12000; M=11500; N="`$/caja 70 unidades"; P=164; Q=70 },
    @{ Row=438; D=44620; H="Sin especificar"; I="Segunda"; J=130; K=9000; L=10000; M=9500; N="`$/caja 100 unidades"; P=95; Q=100 },
    @{ Row=439; D=44606; H="Sin especificar"; I="Primera"; J=120; K=8000; L=9000; M=8500; N="`$/caja 70 unidades"; P=121; Q=70 },
    @{ Row=440; D=44606; H="Sin especificar"; I="Segunda"; J=150; K=7000; L=8000; M=7500; N="`$/caja 100 unidades"; P=75; Q=100 },
    @{ Row=441; D=44567; H="Sin especificar"; I="Primera"; J=120; K=3000; L=3500; M=3250; N="`$/caja 70 unidades"; P=46; Q=70 },
    @{ Row=442; D=44567; H="Sin especificar"; I="Segunda"; J=120; K=2500; L=3000; M=2750; N="`$/caja 100 unidades"; P=28; Q=100 },
    @{ Row=443; D=44960; H="Sin especificar"; I="Primera"; J=150; K=3000; L=4000; M=3400; N="`$/caja 70 unidades"; P=49; Q=70 },
    @{ Row=444; D=44215; H="Sin especificar"; I="Primera"; J=120; K=7000; L=8000; M=7500; N="`$/caja 70 unidades"; P=107; Q=70 },
    @{ Row=445; D=44215; H="Sin especificar"; I="Segunda"; J=120; K=6000; L=7000; M=6500; N="`$/caja 100 unidades"; P=65; Q=100 },
    @{ Row=446; D=44771; H="Sin especificar"; I="Primera"; J=120; K=15000; L=16000; M=15500; N="`$/caja 70 unidades"; P=221; Q=70 },
    @{ Row=447; D=44771; H="Sin especificar"; I="Segunda"; J=170; K=11000; L=12000; M=11500; N="`$/caja 100 unidades"; P=115; Q=100 },
    @{ Row=448; D=44245; H="Sin especificar"; I="Primera"; J=120; K=11000; L=12000; M=11500; N="`$/caja 70 unidades"; P=164; Q=70 },
    @{ Row=449; D=44245; H="Sin especificar"; I="Segunda"; J=120; K=9000; L=10000; M=9500; N="`$/caja 100 unidades"; P=95; Q=100 },
    @{ Row=450; D=44498; H="Sin especificar"; I="Primera"; J=160; K=5000; L=6000; M=5500; N="`$/caja 70 unidades"; P=79; Q=70 },
    @{ Row=451; D=44498; H="Sin especificar"; I="Segunda"; J=160; K=4000; L=5000; M=4500; N="`$/caja 100 unidades"; P=45; Q=100 },
    @{ Row=452; D=44536; H="Sin especificar"; I="Primera"; J=130; K=4000; L=5000; M=4500; N="`$/caja 70 unidades"; P=64; Q=70 },
    @{ Row=453; D=44536; H="Sin especificar"; I="Segunda"; J=120; K=3500; L=4000; M=3750; N="`$/caja 100 unidades"; P=38; Q=100 },
    @{ Row=454; D=44242; H="Sin especificar"; I="Primera"; J=140; K=9000; L=10000; M=9500; N="`$/caja 70 unidades"; P=136; Q=70 },
    @{ Row=455; D=44242; H="Sin especificar"; I="Segunda"; J=130; K=8000; L=9000; M=8500; N="`$/caja 100 unidades"; P=85; Q=100 },
    @{ Row=456; D=44671; H="Sin especificar"; I="Primera"; J=120; K=12000; L=13000; M=12500; N="`$/caja 70 unidades"; P=179; Q=70 },
    @{ Row=457; D=44671; H="Sin especificar"; I="Segunda"; J=150; K=8000; L=9000; M=8500; N="`$/caja 100 unidades"; P=85; Q=100 },
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 4).Value  = $u.D   # D Fecha
    $ws.Cells.Item($r, 8).Value  = $u.H   # H Variedad
    $ws.Cells.Item($r, 9).Value  = $u.I   # I Calidad
    $ws.Cells.Item($r, 10).Value = $u.J   # J Volumen
    $ws.Cells.Item($r, 11).Value = $u.K   # K Precio minimo
    $ws.Cells.Item($r, 12).Value = $u.L   # L Precio maximo
    $ws.Cells.Item($r, 13).Value = $u.M   # M Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $u.N   # N Unidad de comercializacion
    $ws.Cells.Item($r, 16).Value = $u.P   # P Precio $/Kg
    $ws.Cells.Item($r, 17).Value = $u.Q   # Q Kg o Unidades
}

# Append two new rows (458-459) carrying the same constant columns as the rest of
# the block, with D,H,I,J,K,L,M,N,P,Q taken from the tail of the shifted sequence.
$newRows = @(
    @{ Row=458; A=1; B="Agrícola del Norte S.A. de Arica"; C="Arica y Parinacota"; D=44414; E=15; F=100112043; G="Pepino ensalada"; H="Sin especificar"; I="Primera"; J=120; K=13000; L=14000; M=13500; N="`$/caja 70 unidades"; O="Región de Arica y Parinacota"; P=193; Q=70; R="Hortaliza" },
    @{ Row=459; A=1; B="Agrícola del Norte S.A. de Arica"; C="Arica y Parinacota"; D=44414; E=15; F=100112043; G="Pepino ensalada"; H="Sin especificar"; I="Segunda"; J=150; K=8000; L=9000; M=8500; N="`$/caja 100 unidades"; O="Región de Arica y Parinacota"; P=85; Q=100; R="Hortaliza" },
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Cells.Item($r, 1).Value  = $nr.A
    $ws.Cells.Item($r, 2).Value  = $nr.B
    $ws.Cells.Item($r, 3).Value  = $nr.C
    $ws.Cells.Item($r, 4).Value  = $nr.D
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value  = $nr.E
    $ws.Cells.Item($r, 6).Value  = $nr.F
    $ws.Cells.Item($r, 7).Value  = $nr.G
    $ws.Cells.Item($r, 8).Value  = $nr.H
    $ws.Cells.Item($r, 9).Value  = $nr.I
    $ws.Cells.Item($r, 10).Value = $nr.J
    $ws.Cells.Item($r, 11).Value = $nr.K
    $ws.Cells.Item($r, 12).Value = $nr.L
    $ws.Cells.Item($r, 13).Value = $nr.M
    $ws.Cells.Item($r, 14).Value = $nr.N
    $ws.Cells.Item($r, 15).Value = $nr.O
    $ws.Cells.Item($r, 16).Value = $nr.P
    $ws.Cells.Item($r, 17).Value = $nr.Q
    $ws.Cells.Item($r, 18).Value = $nr.R
}